# Weekly update: insert 3 new rows of price data (week of 2022-07-27,
# serial 44769) at the top of the Alcachofa / Comercializadora del Agro de
# Limari price table, pushing the existing rows 213-233 down to 216-236.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 213:233 down by three rows.
$ws.Rows("213:215").Insert()

# Common (constant across every data row in this sheet) column values.
$mercadoId  = 2
$mercado    = "Comercializadora del Agro de Limarí"
$region     = "Coquimbo"
$codreg     = 4
$catId      = 100112013
$categoria  = "Alcachofa"
$clasif     = "Hortaliza"

# --- Row 213: Argentina(o) / Primera ---------------------------------
$ws.Cells.Item(213, 1).Value  = $mercadoId
$ws.Cells.Item(213, 2).Value  = $mercado
$ws.Cells.Item(213, 3).Value  = $region
$ws.Cells.Item(213, 4).Value  = 44769
$ws.Cells.Item(213, 5).Value  = $codreg
$ws.Cells.Item(213, 6).Value  = $catId
$ws.Cells.Item(213, 7).Value  = $categoria
$ws.Cells.Item(213, 8).Value  = "Argentina(o)"
$ws.Cells.Item(213, 9).Value  = "Primera"
$ws.Cells.Item(213, 10).Value = 700
$ws.Cells.Item(213, 11).Value = 9000
$ws.Cells.Item(213, 12).Value = 10000
$ws.Cells.Item(213, 13).Value = 9500
$ws.Cells.Item(213, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(213, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(213, 16).Value = 190
$ws.Cells.Item(213, 17).Value = 50
$ws.Cells.Item(213, 18).Value = $clasif

# --- Row 214: Española / Primera --------------------------------------
$ws.Cells.Item(214, 1).Value  = $mercadoId
$ws.Cells.Item(214, 2).Value  = $mercado
$ws.Cells.Item(214, 3).Value  = $region
$ws.Cells.Item(214, 4).Value  = 44769
$ws.Cells.Item(214, 5).Value  = $codreg
$ws.Cells.Item(214, 6).Value  = $catId
$ws.Cells.Item(214, 7).Value  = $categoria
$ws.Cells.Item(214, 8).Value  = "Española"
$ws.Cells.Item(214, 9).Value  = "Primera"
$ws.Cells.Item(214, 10).Value = 1100
$ws.Cells.Item(214, 11).Value = 11000
$ws.Cells.Item(214, 12).Value = 12000
$ws.Cells.Item(214, 13).Value = 11500
$ws.Cells.Item(214, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(214, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(214, 16).Value = 383
$ws.Cells.Item(214, 17).Value = 30
$ws.Cells.Item(214, 18).Value = $clasif

# --- Row 215: Madrigal / Primera --------------------------------------
$ws.Cells.Item(215, 1).Value  = $mercadoId
$ws.Cells.Item(215, 2).Value  = $mercado
$ws.Cells.Item(215, 3).Value  = $region
$ws.Cells.Item(215, 4).Value  = 44769
$ws.Cells.Item(215, 5).Value  = $codreg
$ws.Cells.Item(215, 6).Value  = $catId
$ws.Cells.Item(215, 7).Value  = $categoria
$ws.Cells.Item(215, 8).Value  = "Madrigal"
$ws.Cells.Item(215, 9).Value  = "Primera"
$ws.Cells.Item(215, 10).Value = 600
$ws.Cells.Item(215, 11).Value = 10000
$ws.Cells.Item(215, 12).Value = 11000
$ws.Cells.Item(215, 13).Value = 10500
$ws.Cells.Item(215, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(215, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(215, 16).Value = 262
$ws.Cells.Item(215, 17).Value = 40
$ws.Cells.Item(215, 18).Value = $clasif
